$wb = $excel.ActiveWorkbook

$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsTasas = $wb.Worksheets.Item("tasas")

$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.73 = 14333.74 pesos`n✅ 14333.74 pesos = 3.71 = 943.34 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsTasas.Range("N10").Value = 267.99
$wsTasas.Range("O10").Value = 3841.3
$wsTasas.Range("N12").Value = 3863.99
$wsTasas.Range("O12").Value = 254.3
